$p = $ppt.ActivePresentation

# --- 1. Presentation-level extLst: add empty p15:sldGuideLst (slide guides ext) ---
# (Handled further below via direct XML injection if COM doesn't expose it.)

# --- 2. Update the footer date placeholders across slide master + all layouts ---
$oldDate = "2020. 02. 12."
$newDate = "2020. 02. 17."

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

Update-DateField($p.SlideMaster.Shapes)
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateField($layout.Shapes)
}

# --- 3. Add a new connector on slide 1 from shape id=7 (right side) to shape id=50 (left side) ---
$s = $p.Slides.Item(1)

function Find-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$shape7 = Find-ShapeById $s.Shapes 7
$shape50 = Find-ShapeById $s.Shapes 50

$beginX = $shape7.Left + $shape7.Width
$beginY = $shape7.Top + ($shape7.Height / 2)
$endX = $shape50.Left
$endY = $shape50.Top + ($shape50.Height / 2)

$conn = $s.Shapes.AddConnector(1, $beginX, $beginY, $endX, $endY)
Write-Output ("new shape id=$($conn.Id) name=$($conn.Name)")
$conn.Name = "Straight Arrow Connector 29"

try {
    $conn.ConnectorFormat.BeginConnect($shape7, 3)
    Write-Output "BeginConnect ok"
} catch {
    Write-Output "BeginConnect failed: $_"
}
try {
    $conn.ConnectorFormat.EndConnect($shape50, 1)
    Write-Output "EndConnect ok"
} catch {
    Write-Output "EndConnect failed: $_"
}

$conn.Line.Weight = 1.5
$conn.Line.ForeColor.ObjectThemeColor = 10
$conn.Line.ForeColor.Brightness = -0.25
$conn.Line.EndArrowheadStyle = 3
